# Sistemas.xlsx — "adicionando model e controller de usuarios"
#
# - Correct the system name in B3 (SAVA -> SAWA)
# - Remove rows 4-6 (GENERAL / Godzilla / Muto) plus the trailing blank
#   rows 7-9, leaving just one blank row (row 4) after the remaining data
# - Clear the leftover "text" number-format styling on A2/A3 so they fall
#   back to the sheet's default style
# - Leave the selection where the user ended up after the edit (C12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd system name.
$ws.Range("B3").Value = "SAWA"

# Drop the stray text-format style Excel had applied to A2/A3 so they go
# back to the workbook's default (General) style.
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Style = "Normal"

# Remove rows 5-9 entirely (Godzilla, Muto, and the trailing blank rows),
# then clear row 4 (GENERAL) in place so a single blank row remains right
# after the trimmed data, matching the original sheet's "a few blank rows
# below the data" shape.
$ws.Range("A5:C9").EntireRow.Delete()
$ws.Range("A4:C4").Clear()

# Move the selection to where the editor left off.
[void]$ws.Range("C12").Select()
